$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3275
$ws.Cells.Item(43, 10).Value = 550
$ws.Cells.Item(43, 12).Value = 550
$ws.Cells.Item(43, 14).Value = -688

$ws.Cells.Item(112, 8).Value = 1944.4445
$ws.Cells.Item(112, 10).Value = 2000
$ws.Cells.Item(112, 12).Value = 6000
$ws.Cells.Item(112, 14).Value = -8216

$ws.Cells.Item(113, 8).Value = 4999.5
$ws.Cells.Item(113, 9).Value = 4999.5
$ws.Cells.Item(113, 11).Value = 4999.5
$ws.Cells.Item(113, 13).Value = -1745.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 631.1818
$ws.Cells.Item(32, 9).Value = 631.1818
$ws.Cells.Item(32, 11).Value = 631.1818
$ws.Cells.Item(32, 13).Value = -344.1818

$ws.Cells.Item(61, 8).Value = 2875.6
$ws.Cells.Item(61, 9).Value = 2875.6
$ws.Cells.Item(61, 11).Value = 2875.6
$ws.Cells.Item(61, 13).Value = -2663.6

$ws.Cells.Item(74, 8).Value = 1897
$ws.Cells.Item(74, 9).Value = 2123.75
$ws.Cells.Item(74, 10).Value = 990
$ws.Cells.Item(74, 11).Value = 2123.75
$ws.Cells.Item(74, 12).Value = 990
$ws.Cells.Item(74, 13).Value = -1249.75
$ws.Cells.Item(74, 14).Value = -2738

$ws.Cells.Item(77, 8).Value = 1897
$ws.Cells.Item(77, 9).Value = 2123.75
$ws.Cells.Item(77, 10).Value = 990
$ws.Cells.Item(77, 11).Value = 10618.75
$ws.Cells.Item(77, 12).Value = 4950
$ws.Cells.Item(77, 13).Value = -6250.75
$ws.Cells.Item(77, 14).Value = -13686

$ws.Cells.Item(132, 8).Value = 1898.375
$ws.Cells.Item(132, 9).Value = 1898.375
$ws.Cells.Item(132, 11).Value = 5695.125
$ws.Cells.Item(132, 13).Value = -3165.125

$ws.Cells.Item(136, 8).Value = 2875.6
$ws.Cells.Item(136, 9).Value = 2875.6
$ws.Cells.Item(136, 11).Value = 8626.799999999999
$ws.Cells.Item(136, 13).Value = -6076.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 957.0769
$ws.Cells.Item(22, 9).Value = 944.4
$ws.Cells.Item(22, 10).Value = 999.3333
$ws.Cells.Item(22, 11).Value = 944.4
$ws.Cells.Item(22, 12).Value = 999.3333
$ws.Cells.Item(22, 13).Value = -771.4
$ws.Cells.Item(22, 14).Value = -1345.3333

$ws.Cells.Item(80, 8).Value = 1672.6666
$ws.Cells.Item(80, 9).Value = 2005.5
$ws.Cells.Item(80, 11).Value = 2005.5
$ws.Cells.Item(80, 13).Value = -1007.5

$ws.Cells.Item(83, 8).Value = 1672.6666
$ws.Cells.Item(83, 9).Value = 2005.5
$ws.Cells.Item(83, 11).Value = 10027.5
$ws.Cells.Item(83, 13).Value = -5035.5

$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 13).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 2823.8572
$ws.Cells.Item(32, 9).Value = 1627.8334
$ws.Cells.Item(32, 11).Value = 1627.8334
$ws.Cells.Item(32, 13).Value = -1311.8334

$ws.Cells.Item(35, 8).Value = 2714.6
$ws.Cells.Item(35, 9).Value = 2714.6
$ws.Cells.Item(35, 11).Value = 2714.6
$ws.Cells.Item(35, 13).Value = -2420.6

$ws.Cells.Item(74, 8).Value = 60987.8
$ws.Cells.Item(74, 10).Value = 60987.8
$ws.Cells.Item(74, 12).Value = 60987.8
$ws.Cells.Item(74, 14).Value = -62735.8

$ws.Cells.Item(77, 8).Value = 60987.8
$ws.Cells.Item(77, 10).Value = 60987.8
$ws.Cells.Item(77, 12).Value = 182963.4
$ws.Cells.Item(77, 14).Value = -191699.4

$ws.Cells.Item(107, 8).Value = 610.9
$ws.Cells.Item(107, 9).Value = 570.4666999999999
$ws.Cells.Item(107, 10).Value = 732.2
$ws.Cells.Item(107, 11).Value = 570.4666999999999
$ws.Cells.Item(107, 12).Value = 732.2
$ws.Cells.Item(107, 13).Value = 1349.5333
$ws.Cells.Item(107, 14).Value = -4572.2

$ws.Cells.Item(141, 8).Value = 344962.66
$ws.Cells.Item(141, 9).Value = 23000
$ws.Cells.Item(141, 10).Value = 988888
$ws.Cells.Item(141, 11).Value = 23000
$ws.Cells.Item(141, 12).Value = 988888
$ws.Cells.Item(141, 13).Value = -17820
$ws.Cells.Item(141, 14).Value = -999248

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1747.9546
$ws.Cells.Item(4, 9).Value = 1595
$ws.Cells.Item(4, 11).Value = 4785
$ws.Cells.Item(4, 13).Value = -4673

$ws.Cells.Item(22, 8).Value = 1996.8
$ws.Cells.Item(22, 9).Value = 1900
$ws.Cells.Item(22, 10).Value = 2001.8948
$ws.Cells.Item(22, 11).Value = 5700
$ws.Cells.Item(22, 12).Value = 6005.6844
$ws.Cells.Item(22, 13).Value = -5531
$ws.Cells.Item(22, 14).Value = -6343.6844

$ws.Cells.Item(26, 8).Value = 1010
$ws.Cells.Item(26, 9).Value = 544
$ws.Cells.Item(26, 10).Value = 1476
$ws.Cells.Item(26, 11).Value = 1632
$ws.Cells.Item(26, 12).Value = 4428
$ws.Cells.Item(26, 13).Value = -1344
$ws.Cells.Item(26, 14).Value = -5004

$ws.Cells.Item(27, 8).Value = 1996.8
$ws.Cells.Item(27, 9).Value = 1900
$ws.Cells.Item(27, 10).Value = 2001.8948
$ws.Cells.Item(27, 11).Value = 5700
$ws.Cells.Item(27, 12).Value = 6005.6844
$ws.Cells.Item(27, 13).Value = -5598
$ws.Cells.Item(27, 14).Value = -6209.6844

$ws.Cells.Item(32, 8).Value = 1990.909
$ws.Cells.Item(32, 9).Value = 91.666664
$ws.Cells.Item(32, 10).Value = 4270
$ws.Cells.Item(32, 11).Value = 274.999992
$ws.Cells.Item(32, 12).Value = 12810
$ws.Cells.Item(32, 13).Value = 8.00000799999998
$ws.Cells.Item(32, 14).Value = -13376

$ws.Cells.Item(33, 8).Value = 39.75
$ws.Cells.Item(33, 9).Value = 40
$ws.Cells.Item(33, 10).Value = 39.5
$ws.Cells.Item(33, 11).Value = 240
$ws.Cells.Item(33, 12).Value = 237
$ws.Cells.Item(33, 13).Value = 43
$ws.Cells.Item(33, 14).Value = -803

$ws.Cells.Item(34, 8).Value = 34640.594
$ws.Cells.Item(34, 9).Value = 91.59999999999999
$ws.Cells.Item(34, 10).Value = 42492.637
$ws.Cells.Item(34, 11).Value = 274.8
$ws.Cells.Item(34, 12).Value = 127477.911
$ws.Cells.Item(34, 13).Value = -190.8
$ws.Cells.Item(34, 14).Value = -127645.911

$ws.Cells.Item(38, 8).Value = 174
$ws.Cells.Item(38, 9).Value = 155
$ws.Cells.Item(38, 11).Value = 465
$ws.Cells.Item(38, 13).Value = -118

$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 13).ClearContents()

$ws.Cells.Item(40, 8).Value = 133.53847
$ws.Cells.Item(40, 9).Value = 65.666664
$ws.Cells.Item(40, 10).Value = 286.25
$ws.Cells.Item(40, 11).Value = 262.666656
$ws.Cells.Item(40, 12).Value = 1145
$ws.Cells.Item(40, 13).Value = -193.666656
$ws.Cells.Item(40, 14).Value = -1283

$ws.Cells.Item(44, 8).Value = 931.23334
$ws.Cells.Item(44, 9).Value = 463
$ws.Cells.Item(44, 10).Value = 947.37933
$ws.Cells.Item(44, 11).Value = 1389
$ws.Cells.Item(44, 12).Value = 2842.13799
$ws.Cells.Item(44, 13).Value = -991
$ws.Cells.Item(44, 14).Value = -3638.13799

$ws.Cells.Item(46, 8).Value = 2029.2307
$ws.Cells.Item(46, 9).Value = 398.83334
$ws.Cells.Item(46, 10).Value = 3426.7144
$ws.Cells.Item(46, 11).Value = 1196.50002
$ws.Cells.Item(46, 12).Value = 10280.1432
$ws.Cells.Item(46, 13).Value = -1105.50002
$ws.Cells.Item(46, 14).Value = -10462.1432

$ws.Cells.Item(64, 8).Value = 2690
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2690
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 8070
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -8610

$ws.Cells.Item(67, 8).Value = 2690
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 2690
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 8070
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -9942

$ws.Cells.Item(86, 8).Value = 558.6667
$ws.Cells.Item(86, 9).Value = 251
$ws.Cells.Item(86, 10).Value = 712.5
$ws.Cells.Item(86, 11).Value = 753
$ws.Cells.Item(86, 12).Value = 2137.5
$ws.Cells.Item(86, 13).Value = 433
$ws.Cells.Item(86, 14).Value = -4509.5

$ws.Cells.Item(89, 8).Value = 558.6667
$ws.Cells.Item(89, 9).Value = 251
$ws.Cells.Item(89, 10).Value = 712.5
$ws.Cells.Item(89, 11).Value = 2259
$ws.Cells.Item(89, 12).Value = 6412.5
$ws.Cells.Item(89, 13).Value = 3669
$ws.Cells.Item(89, 14).Value = -18268.5

$ws.Cells.Item(120, 8).Value = 8944.333000000001
$ws.Cells.Item(120, 9).Value = 5900
$ws.Cells.Item(120, 10).Value = 15033
$ws.Cells.Item(120, 11).Value = 17700
$ws.Cells.Item(120, 12).Value = 45099
$ws.Cells.Item(120, 13).Value = -12862
$ws.Cells.Item(120, 14).Value = -54775

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 19000
$ws.Cells.Item(7, 9).Value = 19000
$ws.Cells.Item(7, 11).Value = 19000
$ws.Cells.Item(7, 13).Value = -18888

$ws.Cells.Item(22, 8).Value = 4000
$ws.Cells.Item(22, 9).Value = 5000
$ws.Cells.Item(22, 10).Value = 3000
$ws.Cells.Item(22, 11).Value = 5000
$ws.Cells.Item(22, 12).Value = 3000
$ws.Cells.Item(22, 13).Value = -4705
$ws.Cells.Item(22, 14).Value = -3590

$ws.Cells.Item(27, 8).Value = 4000
$ws.Cells.Item(27, 9).Value = 5000
$ws.Cells.Item(27, 10).Value = 3000
$ws.Cells.Item(27, 11).Value = 5000
$ws.Cells.Item(27, 12).Value = 3000
$ws.Cells.Item(27, 13).Value = -4893
$ws.Cells.Item(27, 14).Value = -3214

$ws.Cells.Item(64, 8).Value = 87499.5
$ws.Cells.Item(64, 10).Value = 87499.5
$ws.Cells.Item(64, 12).Value = 87499.5
$ws.Cells.Item(64, 14).Value = -87949.5

$ws.Cells.Item(67, 8).Value = 87499.5
$ws.Cells.Item(67, 10).Value = 87499.5
$ws.Cells.Item(67, 12).Value = 87499.5
$ws.Cells.Item(67, 14).Value = -89059.5

$ws.Cells.Item(122, 8).Value = 3355.889
$ws.Cells.Item(122, 9).Value = 2939.6
$ws.Cells.Item(122, 11).Value = 8818.799999999999
$ws.Cells.Item(122, 13).Value = -6368.799999999999

$ws.Cells.Item(126, 8).Value = 19000
$ws.Cells.Item(126, 9).Value = 19000
$ws.Cells.Item(126, 11).Value = 57000
$ws.Cells.Item(126, 13).Value = -54530

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1509.8334
$ws.Cells.Item(81, 9).Value = 1423.2222
$ws.Cells.Item(81, 10).Value = 1769.6666
$ws.Cells.Item(81, 11).Value = 2846.4444
$ws.Cells.Item(81, 12).Value = 3539.3332
$ws.Cells.Item(81, 13).Value = -1785.4444
$ws.Cells.Item(81, 14).Value = -5661.3332

$ws.Cells.Item(84, 8).Value = 1509.8334
$ws.Cells.Item(84, 9).Value = 1423.2222
$ws.Cells.Item(84, 10).Value = 1769.6666
$ws.Cells.Item(84, 11).Value = 14232.222
$ws.Cells.Item(84, 12).Value = 17696.666
$ws.Cells.Item(84, 13).Value = -8928.222
$ws.Cells.Item(84, 14).Value = -28304.666

$ws.Cells.Item(113, 8).Value = 222.23077
$ws.Cells.Item(113, 9).Value = 144.3
$ws.Cells.Item(113, 11).Value = 432.9
$ws.Cells.Item(113, 13).Value = 1737.1

$ws.Cells.Item(136, 8).Value = 2377
$ws.Cells.Item(136, 9).Value = 2377
$ws.Cells.Item(136, 11).Value = 7131
$ws.Cells.Item(136, 13).Value = -4581
